$d = $word.ActiveDocument

# Collapse a range to the very end of the document body (just before sectPr)
$endRange = $d.Content
$endRange.Collapse(0)

# Build the raw OOXML for the five new paragraphs that follow the last
# existing paragraph ("...甚至可以不定义构造函数。") and precede <w:sectPr>:
#   1. an empty paragraph
#   2. an empty paragraph
#   3. a paragraph with explicit formatting holding the copyright text
#   4. a paragraph with the same formatting but no text
#   5. a paragraph carrying only an eastAsia font hint on the mark
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = '<w:p ' + $ns + '/>'
$xml += '<w:p ' + $ns + '/>'
$xml += '<w:p ' + $ns + '>' +
          '<w:pPr>' +
            '<w:widowControl/>' +
            '<w:jc w:val="left"/>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体"/>' +
              '<w:kern w:val="0"/>' +
              '<w:szCs w:val="21"/>' +
            '</w:rPr>' +
          '</w:pPr>' +
          '<w:r>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体" w:hint="eastAsia"/>' +
              '<w:kern w:val="0"/>' +
              '<w:szCs w:val="21"/>' +
            '</w:rPr>' +
            '<w:t>Copyright ©2021-2099 HaoyangZheng. All rights reserved</w:t>' +
          '</w:r>' +
        '</w:p>'
$xml += '<w:p ' + $ns + '>' +
          '<w:pPr>' +
            '<w:widowControl/>' +
            '<w:jc w:val="left"/>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体" w:hint="eastAsia"/>' +
              '<w:kern w:val="0"/>' +
              '<w:szCs w:val="21"/>' +
            '</w:rPr>' +
          '</w:pPr>' +
        '</w:p>'
$xml += '<w:p ' + $ns + '>' +
          '<w:pPr>' +
            '<w:rPr>' +
              '<w:rFonts w:hint="eastAsia"/>' +
            '</w:rPr>' +
          '</w:pPr>' +
        '</w:p>'

[void]$endRange.InsertXML($xml)
